$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 113 (pushes existing rows 113-126 down to 114-127)
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with the new weekly data point
$ws.Cells.Item(113, 1).Value = 1
$ws.Cells.Item(113, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(113, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(113, 4).Value = 45166
$ws.Cells.Item(113, 5).Value = 15
$ws.Cells.Item(113, 6).Value = 100112040
$ws.Cells.Item(113, 7).Value = "Cilantro"
$ws.Cells.Item(113, 8).Value = "Sin especificar"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 300
$ws.Cells.Item(113, 11).Value = 1500
$ws.Cells.Item(113, 12).Value = 2000
$ws.Cells.Item(113, 13).Value = 1750
$ws.Cells.Item(113, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(113, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(113, 16).Value = 875
$ws.Cells.Item(113, 17).Value = 2
$ws.Cells.Item(113, 18).Value = "Hortaliza"

# Keep date formatting consistent with the rest of column D
$ws.Cells.Item(113, 4).NumberFormat = $ws.Cells.Item(114, 4).NumberFormat
